$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.570.63'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '3.514.99'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''600.85'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = '''142.62'
$ws.Range('D7').Value = '3.515.06'
$ws.Range('E7').Value = '  -2.07%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '''0.516'
$ws.Range('E9').Value = '  +5.19%  '
$ws.Range('E10').Value = '  -2.88%  '
$ws.Range('D11').Value = '''7.82'
$ws.Range('E11').Value = '  -2.27%  '
$ws.Range('D12').Value = '''0.404'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('D13').Value = '4.115.62'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = '''0.0000196'
$ws.Range('E14').Value = '  -6.24%  '
$ws.Range('D15').Value = '''28.35'
$ws.Range('E15').Value = '  -5.86%  '
$ws.Range('D16').Value = '3.506.83'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '65.492.79'
$ws.Range('E18').Value = '  -1.86%  '
$ws.Range('D19').Value = '''10.85'
$ws.Range('E19').Value = '  -5.37%  '
$ws.Range('D20').Value = '''6.16'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('D21').Value = '''14.41'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('D22').Value = '''417.96'
$ws.Range('E22').Value = '  -3.22%  '
$ws.Range('D23').Value = '''0.595'
$ws.Range('E23').Value = '  -4.79%  '
$ws.Range('D24').Value = '''76.99'
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('D25').Value = '3.654.45'
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '''0.0000114'
$ws.Range('E27').Value = '  -5.35%  '
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''7.73'
$ws.Range('E29').Value = '  -5.74%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''8.88'
$ws.Range('E30').Value = '  -4.73%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = '3.522.82'
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('D34').Value = '''24.19'
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '''1.34'
$ws.Range('E36').Value = '  -7.96%  '
$ws.Range('D37').Value = '''7.53'
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').Value = '''174.12'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').Value = '''5.23'
$ws.Range('E39').Value = '  -7.16%  '
$ws.Range('D40').Value = '''1.57'
$ws.Range('E40').Value = '  -8.73%  '
$ws.Range('D41').Value = '''0.0810'
$ws.Range('E41').Value = '  -5.56%  '
$ws.Range('D42').Value = '''5.00'
$ws.Range('E42').Value = '  -4.64%  '
$ws.Range('D43').Value = '''0.855'
$ws.Range('E43').Value = '  -4.64%  '
$ws.Range('D44').Value = '''45.23'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('D45').Value = '''1.76'
$ws.Range('E45').Value = '  -8.05%  '
$ws.Range('D46').Value = '''0.999'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  -8.40%  '
$ws.Range('D48').Value = '''7.03'
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').Value = '''23.12'
$ws.Range('E49').Value = '  -3.65%  '
$ws.Range('E50').Value = '  -8.86%  '
$ws.Range('D51').Value = '''0.902'
$ws.Range('E51').Value = '  -5.13%  '
